$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values for existing rows 2-34, and add new rows 35-42 (A: date serial, B: count, C: shared "M")
$ws.Cells.Item(2, 1).Value = 43496
$ws.Cells.Item(2, 2).Value = 1658
$ws.Cells.Item(2, 3).Value = "M"

$ws.Cells.Item(3, 1).Value = 43524
$ws.Cells.Item(3, 2).Value = 1683
$ws.Cells.Item(3, 3).Value = "M"

$ws.Cells.Item(4, 1).Value = 43555
$ws.Cells.Item(4, 2).Value = 1686
$ws.Cells.Item(4, 3).Value = "M"

$ws.Cells.Item(5, 1).Value = 43585
$ws.Cells.Item(5, 2).Value = 1728
$ws.Cells.Item(5, 3).Value = "M"

$ws.Cells.Item(6, 1).Value = 43616
$ws.Cells.Item(6, 2).Value = 1679
$ws.Cells.Item(6, 3).Value = "M"

$ws.Cells.Item(7, 1).Value = 43646
$ws.Cells.Item(7, 2).Value = 1643
$ws.Cells.Item(7, 3).Value = "M"

$ws.Cells.Item(8, 1).Value = 43677
$ws.Cells.Item(8, 2).Value = 1640
$ws.Cells.Item(8, 3).Value = "M"

$ws.Cells.Item(9, 1).Value = 43708
$ws.Cells.Item(9, 2).Value = 1623
$ws.Cells.Item(9, 3).Value = "M"

$ws.Cells.Item(10, 1).Value = 43738
$ws.Cells.Item(10, 2).Value = 1640
$ws.Cells.Item(10, 3).Value = "M"

$ws.Cells.Item(11, 1).Value = 43769
$ws.Cells.Item(11, 2).Value = 1626
$ws.Cells.Item(11, 3).Value = "M"

$ws.Cells.Item(12, 1).Value = 43799
$ws.Cells.Item(12, 2).Value = 1622
$ws.Cells.Item(12, 3).Value = "M"

$ws.Cells.Item(13, 1).Value = 43830
$ws.Cells.Item(13, 2).Value = 1560
$ws.Cells.Item(13, 3).Value = "M"

$ws.Cells.Item(14, 1).Value = 43861
$ws.Cells.Item(14, 2).Value = 1575
$ws.Cells.Item(14, 3).Value = "M"

$ws.Cells.Item(15, 1).Value = 43890
$ws.Cells.Item(15, 2).Value = 1602
$ws.Cells.Item(15, 3).Value = "M"

$ws.Cells.Item(16, 1).Value = 43921
$ws.Cells.Item(16, 2).Value = 1660
$ws.Cells.Item(16, 3).Value = "M"

$ws.Cells.Item(17, 1).Value = 43951
$ws.Cells.Item(17, 2).Value = 1528
$ws.Cells.Item(17, 3).Value = "M"

$ws.Cells.Item(18, 1).Value = 43982
$ws.Cells.Item(18, 2).Value = 1475
$ws.Cells.Item(18, 3).Value = "M"

$ws.Cells.Item(19, 1).Value = 44012
$ws.Cells.Item(19, 2).Value = 1495
$ws.Cells.Item(19, 3).Value = "M"

$ws.Cells.Item(20, 1).Value = 44043
$ws.Cells.Item(20, 2).Value = 1509
$ws.Cells.Item(20, 3).Value = "M"

$ws.Cells.Item(21, 1).Value = 44074
$ws.Cells.Item(21, 2).Value = 1528
$ws.Cells.Item(21, 3).Value = "M"

$ws.Cells.Item(22, 1).Value = 44104
$ws.Cells.Item(22, 2).Value = 1509
$ws.Cells.Item(22, 3).Value = "M"

$ws.Cells.Item(23, 1).Value = 44135
$ws.Cells.Item(23, 2).Value = 1547
$ws.Cells.Item(23, 3).Value = "M"

$ws.Cells.Item(24, 1).Value = 44165
$ws.Cells.Item(24, 2).Value = 1531
$ws.Cells.Item(24, 3).Value = "M"

$ws.Cells.Item(25, 1).Value = 44196
$ws.Cells.Item(25, 2).Value = 1513
$ws.Cells.Item(25, 3).Value = "M"

$ws.Cells.Item(26, 1).Value = 44227
$ws.Cells.Item(26, 2).Value = 1456
$ws.Cells.Item(26, 3).Value = "M"

$ws.Cells.Item(27, 1).Value = 44255
$ws.Cells.Item(27, 2).Value = 1469
$ws.Cells.Item(27, 3).Value = "M"

$ws.Cells.Item(28, 1).Value = 44286
$ws.Cells.Item(28, 2).Value = 1501
$ws.Cells.Item(28, 3).Value = "M"

$ws.Cells.Item(29, 1).Value = 44316
$ws.Cells.Item(29, 2).Value = 1532
$ws.Cells.Item(29, 3).Value = "M"

$ws.Cells.Item(30, 1).Value = 44347
$ws.Cells.Item(30, 2).Value = 1554
$ws.Cells.Item(30, 3).Value = "M"

$ws.Cells.Item(31, 1).Value = 44377
$ws.Cells.Item(31, 2).Value = 1542
$ws.Cells.Item(31, 3).Value = "M"

$ws.Cells.Item(32, 1).Value = 44408
$ws.Cells.Item(32, 2).Value = 1540
$ws.Cells.Item(32, 3).Value = "M"

$ws.Cells.Item(33, 1).Value = 44439
$ws.Cells.Item(33, 2).Value = 1503
$ws.Cells.Item(33, 3).Value = "M"

$ws.Cells.Item(34, 1).Value = 44469
$ws.Cells.Item(34, 2).Value = 1504
$ws.Cells.Item(34, 3).Value = "M"

$ws.Cells.Item(35, 1).Value = 44500
$ws.Cells.Item(35, 2).Value = 1495
$ws.Cells.Item(35, 3).Value = "M"
$ws.Cells.Item(35, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(36, 1).Value = 44530
$ws.Cells.Item(36, 2).Value = 1482
$ws.Cells.Item(36, 3).Value = "M"
$ws.Cells.Item(36, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(37, 1).Value = 44561
$ws.Cells.Item(37, 2).Value = 1433
$ws.Cells.Item(37, 3).Value = "M"
$ws.Cells.Item(37, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(38, 1).Value = 44592
$ws.Cells.Item(38, 2).Value = 1442
$ws.Cells.Item(38, 3).Value = "M"
$ws.Cells.Item(38, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(39, 1).Value = 44620
$ws.Cells.Item(39, 2).Value = 1488
$ws.Cells.Item(39, 3).Value = "M"
$ws.Cells.Item(39, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(40, 1).Value = 44651
$ws.Cells.Item(40, 2).Value = 1494
$ws.Cells.Item(40, 3).Value = "M"
$ws.Cells.Item(40, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(41, 1).Value = 44681
$ws.Cells.Item(41, 2).Value = 1500
$ws.Cells.Item(41, 3).Value = "M"
$ws.Cells.Item(41, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

$ws.Cells.Item(42, 1).Value = 44712
$ws.Cells.Item(42, 2).Value = 1511
$ws.Cells.Item(42, 3).Value = "M"
$ws.Cells.Item(42, 1).NumberFormat = $ws.Cells.Item(34, 1).NumberFormat

